$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "JSU(-0.9478040347888075, 1.1600942729923163, 0.38058756717037434, 2.6196203984555417)"
$ws.Range("C2").Value = "NIG(1.6378198977583236, 1.2613032702141227, 3.3501997243659543, 6.048688549293481)"
$ws.Range("D2").Value = "JSU(-0.8153266048216445, 1.2247826064387488, 1.4796470569435838, 3.2665743593335295)"
$ws.Range("E2").Value = "EXN(1.8507580504911991, 3.440136965907115, 3.804158205796213)"
